$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buffer Days Plan")
$ws.Range("E80").Value = "x"
$ws.Range("E81").Value = "x"
$ws.Range("E82").Value = "x"
$ws.Range("E83").Value = "x"
$ws.Range("E84").Value = "x"
$ws.Range("E86").Value = "x"
$ws.Range("E87").Value = "x"
$ws.Range("E88").Value = "x"
